$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = '2025-12-19 01:24:39'

$arr = New-Object 'object[,]' 15,8

# Row 2: 【急募】AI医療系請求IOSアプリ開発のフリーランス募集
$arr[0,0] = $timestamp
$arr[0,1] = '【急募】AI医療系請求IOSアプリ開発のフリーランス募集'
$arr[0,2] = 'システム開発'
$arr[0,3] = '1,000,000 円 ~ 3,000,000 円 / 固定'
$arr[0,4] = '期限情報なし'
$arr[0,5] = 'https://www.lancers.jp/work/detail/5456942'
$arr[0,6] = 385
$arr[0,7] = '🔥AI,Ai ◆開発 ◇アプリ'

# Row 3: EC×AIプロダクト/業務改善リード
$arr[1,0] = $timestamp
$arr[1,1] = 'EC×AIプロダクト/業務改善リード'
$arr[1,2] = 'システム開発'
$arr[1,3] = '200,000 円 ~ 300,000 円 / 固定'
$arr[1,4] = '期限情報なし'
$arr[1,5] = 'https://www.lancers.jp/work/detail/5450024'
$arr[1,6] = 338
$arr[1,7] = '🔥AI,Ai ◇業務改善'

# Row 4: 初回 既存システムのRuby、Ruby on Railsバー
$arr[2,0] = $timestamp
$arr[2,1] = '初回 既存システムのRuby、Ruby on Railsバージョンアップ及び追加改修'
$arr[2,2] = 'システム開発'
$arr[2,3] = '50,000 円 ~ 100,000 円 / 固定'
$arr[2,4] = '期限情報なし'
$arr[2,5] = 'https://www.lancers.jp/work/detail/5456434'
$arr[2,6] = 318
$arr[2,7] = '🔥AI,Ai'

# Row 5: 製造業のR&D支援!「プロセスデータ解析」「音響異常検知」の
$arr[3,0] = $timestamp
$arr[3,1] = '製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集'
$arr[3,2] = 'システム開発'
$arr[3,3] = '200,000 円 ~ 300,000 円 / 固定'
$arr[3,4] = '期限情報なし'
$arr[3,5] = 'https://www.lancers.jp/work/detail/5439165'
$arr[3,6] = 303
$arr[3,7] = '🔥AI,Ai'

# Row 6: 【Zapier設定のみ!作業時間~2時間】スプレッドシート・
$arr[4,0] = $timestamp
$arr[4,1] = '【Zapier設定のみ!作業時間~2時間】スプレッドシート・Gドライブ自動化構築(設計済)'
$arr[4,2] = 'システム開発'
$arr[4,3] = '10,000 円 ~ 20,000 円 / 固定'
$arr[4,4] = '期限情報なし'
$arr[4,5] = 'https://www.lancers.jp/work/detail/5456066'
$arr[4,6] = 255
$arr[4,7] = '🔥API ◆自動化'

# Row 7: 【スマホアプリ開発】 音声データ推定アプリの依頼
$arr[5,0] = $timestamp
$arr[5,1] = '【スマホアプリ開発】 音声データ推定アプリの依頼'
$arr[5,2] = 'システム開発'
$arr[5,3] = '500,000 円 ~ 1,000,000 円 / 固定'
$arr[5,4] = '期限情報なし'
$arr[5,5] = 'https://www.lancers.jp/work/detail/5456360'
$arr[5,6] = 175
$arr[5,7] = '★スマホアプリ ◆開発 ◇アプリ'

# Row 8: 【急募】多店舗パーソナルジム向け予約・顧客管理システム開発
$arr[6,0] = $timestamp
$arr[6,1] = '【急募】多店舗パーソナルジム向け予約・顧客管理システム開発'
$arr[6,2] = 'システム開発'
$arr[6,3] = '1,000,000 円 ~ 3,000,000 円 / 固定'
$arr[6,4] = '期限情報なし'
$arr[6,5] = 'https://www.lancers.jp/work/detail/5456461'
$arr[6,6] = 160
$arr[6,7] = '◆開発,システム開発 ◇管理'

# Row 9: 【急募】大規模データ収集自動化(スクレイピング・DB連携・エ
$arr[7,0] = $timestamp
$arr[7,1] = '【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件'
$arr[7,2] = 'システム開発'
$arr[7,3] = '100,000 円 ~ 200,000 円 / 固定'
$arr[7,4] = '期限情報なし'
$arr[7,5] = 'https://www.lancers.jp/work/detail/5456658'
$arr[7,6] = 158
$arr[7,7] = '◆自動化,スクレイピング ◇管理'

# Row 10: 【急募】飲食店予約サイトの制作と将来的なアプリ化(アプリ化の
$arr[8,0] = $timestamp
$arr[8,1] = '【急募】飲食店予約サイトの制作と将来的なアプリ化(アプリ化の際は別契約)'
$arr[8,2] = 'システム開発'
$arr[8,3] = '1,000,000 円 ~ 3,000,000 円 / 固定'
$arr[8,4] = '期限情報なし'
$arr[8,5] = 'https://www.lancers.jp/work/detail/5457089'
$arr[8,6] = 70
$arr[8,7] = '◇アプリ'

# Row 11: 【完全在宅/時給1,400円】IT・業務効率化経験を活かせる
$arr[9,0] = $timestamp
$arr[9,1] = '【完全在宅/時給1,400円】IT・業務効率化経験を活かせる!社内エンジニア兼総務スタッフを募集!'
$arr[9,2] = 'システム開発'
$arr[9,3] = '1,000 ~ 5,000 円 / 固定'
$arr[9,4] = '期限情報なし'
$arr[9,5] = 'https://www.lancers.jp/work/detail/5456452'
$arr[9,6] = 70
$arr[9,7] = '◆効率化'

# Row 12: 【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加
$arr[10,0] = $timestamp
$arr[10,1] = '【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集'
$arr[10,2] = 'システム開発'
$arr[10,3] = '20,000 円 ~ 50,000 円 / 固定'
$arr[10,4] = '期限情報なし'
$arr[10,5] = 'https://www.lancers.jp/work/detail/5457026'
$arr[10,6] = 68
$arr[10,7] = '◆ツール'

# Row 13: 【急募】PHPによる申請サイト構築支援!
$arr[11,0] = $timestamp
$arr[11,1] = '【急募】PHPによる申請サイト構築支援!'
$arr[11,2] = 'システム開発'
$arr[11,3] = '50,000 円 ~ 100,000 円 / 固定'
$arr[11,4] = '期限情報なし'
$arr[11,5] = 'https://www.lancers.jp/work/detail/5457023'
$arr[11,6] = 58
$arr[11,7] = '◇サイト ○PHP'

# Row 14: 【急募】データ活用インフラ要件整理のコンサルタント募集(1人
$arr[12,0] = $timestamp
$arr[12,1] = '【急募】データ活用インフラ要件整理のコンサルタント募集(1人月/月)'
$arr[12,2] = 'システム開発'
$arr[12,3] = '1,000,000 円 ~ 3,000,000 円 / 固定'
$arr[12,4] = '期限情報なし'
$arr[12,5] = 'https://www.lancers.jp/work/detail/5456545'
$arr[12,6] = 55
$arr[12,7] = '◆コンサル'

# Row 15: 【急募】LINEシステム構築・保守運用のプロフェッショナルを
$arr[13,0] = $timestamp
$arr[13,1] = '【急募】LINEシステム構築・保守運用のプロフェッショナルを求む!'
$arr[13,2] = 'システム開発'
$arr[13,3] = '50,000 円 ~ 100,000 円 / 固定'
$arr[13,4] = '期限情報なし'
$arr[13,5] = 'https://www.lancers.jp/work/detail/5456063'
$arr[13,6] = 33
$arr[13,7] = $null

# Row 16: 【準委任】音声データ収集プロジェクト/PM・ディレクター募集
$arr[14,0] = $timestamp
$arr[14,1] = '【準委任】音声データ収集プロジェクト/PM・ディレクター募集'
$arr[14,2] = 'システム開発'
$arr[14,3] = '300,000 円 ~ 500,000 円 / 固定'
$arr[14,4] = '期限情報なし'
$arr[14,5] = 'https://www.lancers.jp/work/detail/5456449'
$arr[14,6] = 25
$arr[14,7] = $null

$ws.Range("A2:H16").Value = $arr

# Rebuild hyperlinks for F2:F16 from scratch to match final expected state
$ws.Hyperlinks.Delete()

$urls = @(
    'https://www.lancers.jp/work/detail/5456942',
    'https://www.lancers.jp/work/detail/5450024',
    'https://www.lancers.jp/work/detail/5456434',
    'https://www.lancers.jp/work/detail/5439165',
    'https://www.lancers.jp/work/detail/5456066',
    'https://www.lancers.jp/work/detail/5456360',
    'https://www.lancers.jp/work/detail/5456461',
    'https://www.lancers.jp/work/detail/5456658',
    'https://www.lancers.jp/work/detail/5457089',
    'https://www.lancers.jp/work/detail/5456452',
    'https://www.lancers.jp/work/detail/5457026',
    'https://www.lancers.jp/work/detail/5457023',
    'https://www.lancers.jp/work/detail/5456545',
    'https://www.lancers.jp/work/detail/5456063',
    'https://www.lancers.jp/work/detail/5456449'
)

for ($i = 0; $i -lt 15; $i++) {
    $rowNum = $i + 2
    $cell = $ws.Cells.Item($rowNum, 6)
    $url = $urls[$i]
    $ws.Hyperlinks.Add($cell, $url, "", "", $url) | Out-Null
    $cell.Style = "Hyperlink"
}
